# fix NPC HP error
# Insert a new "MAXHP" column before the existing "MAXMP" column (column G),
# shifting all subsequent columns one to the right, and fix the HP values
# that were erroneously left at 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G (old column G "MAXMP" and everything after shifts right by one)
$ws.Columns("G").Insert()

# Keep the same column width as column F (14) for the newly inserted column G
$ws.Columns("G").ColumnWidth = $ws.Columns("F").ColumnWidth

# Header for the new column
$ws.Range("G1").Value2 = "MAXHP"

# Fill in MAXHP (new col G) and correct MAXMP (now shifted to col H) values,
# both should equal the NPC's SalePrice-derived HP/MP value (col F) instead of 0.
$lastRow = 6
for ($r = 2; $r -le $lastRow; $r++) {
    $salePrice = $ws.Cells.Item($r, 6).Value2
    $ws.Cells.Item($r, 7).Value2 = $salePrice
    $ws.Cells.Item($r, 8).Value2 = $salePrice
}

# Update the active selection/view to match the saved state
$ws.Range("I8").Select() | Out-Null
